# Apply the UAT_Reporting_Data_Linkage_Requirements edit:
#  - Renumber the STUDIES section items from "2.x" to "1.x"
#  - Remove the stray "2" numeric marker in A5
#  - Add a new "Study-level Consent Details Report" section heading in B19
#    (bold dark-green text, matching the other section headers)
#  - Update the sheet view (scrolled position / active selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the stray "2" value that used to sit above the STUDIES heading
$ws.Range("A5").ClearContents()

# Renumber "2.1".."2.12" -> "1.1".."1.12"
$ws.Range("A6").Value2  = "1.1"
$ws.Range("A7").Value2  = "1.2"
$ws.Range("A8").Value2  = "1.3"
$ws.Range("A9").Value2  = "1.4"
$ws.Range("A10").Value2 = "1.5"
$ws.Range("A11").Value2 = "1.6"
$ws.Range("A12").Value2 = "1.7"
$ws.Range("A13").Value2 = "1.8"
$ws.Range("A14").Value2 = "1.9"
$ws.Range("A15").Value2 = "1.10"
$ws.Range("A16").Value2 = "1.11"
$ws.Range("A17").Value2 = "1.12"

# Add the new section heading into B19, styled like the other section
# headers (same formatting as B18 but bold, dark green text)
$ws.Range("B19").Value2 = "Study-level Consent Details Report"
$ws.Range("B19").Font.Name = "Calibri"
$ws.Range("B19").Font.Size = 11
$ws.Range("B19").Font.Bold = $true
$ws.Range("B19").Font.Color = 24832
$ws.Range("B19").HorizontalAlignment = -4131
$ws.Range("B19").VerticalAlignment = -4160
$ws.Range("B19").WrapText = $true
$ws.Range("B19").Locked = $true

# Update the view so it is scrolled to/selecting the new row
$ws.Range("B19").Select() | Out-Null
